# Update the "Fruta, Vega Monumental Concepción - Caqui" sheet with the latest
# weekly price data (Fruta / hortaliza, semanal).
#
# The underlying rows (2-19) keep the same fixed attributes (Mercado, Región,
# Codreg, Tipo, Producto, Categoría, Variedad) but the weekly observations in
# columns D (Fecha), L (Calidad), M (Volumen), N (Precio mínimo),
# O (Precio máximo), P (Precio promedio ponderado), Q (Unidad de
# comercialización), R (Origen), S (Precio $/Kg) and T (Kg / unidad) are
# refreshed/reordered for the new week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45097, 'Primera', 100, 18000, 20000, 19000, '$/caja 18 kilos granel', 'Región del Maule', 1056, 18),
    @(3, 45077, 'Primera', 140, 12000, 14000, 12857, '$/caja 12 kilos granel', 'Región de O''Higgins', 12857, 1),
    @(4, 45077, 'Segunda', 80, 11000, 11000, 11000, '$/caja 12 kilos granel', 'Región de O''Higgins', 11000, 1),
    @(5, 45084, 'Primera', 100, 17000, 18000, 17500, '$/caja 18 kilos granel', 'Región del Maule', 972, 18),
    @(6, 45100, 'Primera', 60, 18000, 18000, 18000, '$/caja 18 kilos granel', 'Región de O''Higgins', 1000, 18),
    @(7, 44719, 'Primera', 50, 14000, 15000, 14400, '$/caja 18 kilos granel', 'Región del Maule', 800, 18),
    @(8, 45091, 'Primera', 220, 18000, 19000, 18455, '$/caja 18 kilos granel', 'Provincia de Curicó', 1025, 18),
    @(9, 45091, 'Segunda', 150, 15000, 15000, 15000, '$/caja 18 kilos granel', 'Provincia de Curicó', 833, 18),
    @(10, 44707, 'Primera', 60, 12000, 13000, 12500, '$/caja 12 kilos empedrada', 'Provincia de Curicó', 1042, 12),
    @(11, 44334, 'Primera', 100, 11000, 12000, 11500, '$/caja 12 kilos granel', 'Región de O''Higgins', 11500, 1),
    @(12, 45092, 'Primera', 140, 18000, 19000, 18429, '$/caja 18 kilos granel', 'Provincia de Curicó', 1024, 18),
    @(13, 44708, 'Primera', 70, 12000, 13000, 12571, '$/caja 12 kilos empedrada', 'Provincia de Curicó', 1048, 12),
    @(14, 45093, 'Primera', 140, 17000, 18000, 17429, '$/caja 18 kilos granel', 'Provincia de Curicó', 968, 18),
    @(15, 45090, 'Primera', 150, 17000, 18000, 17533, '$/caja 18 kilos granel', 'Región del Maule', 974, 18),
    @(16, 45090, 'Segunda', 130, 14000, 15000, 14462, '$/caja 18 kilos granel', 'Región del Maule', 803, 18),
    @(17, 44714, 'Primera', 100, 14000, 15000, 14500, '$/caja 18 kilos granel', 'Región de O''Higgins', 806, 18),
    @(18, 44330, 'Primera', 100, 15000, 16000, 15500, '$/caja 18 kilos granel', 'Provincia de Curicó', 861, 18),
    @(19, 44742, 'Segunda', 100, 14000, 15000, 14500, '$/caja 18 kilos granel', 'Región de O''Higgins', 806, 18),
)

foreach ($row in $data) {
    $r    = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D - Fecha
    $ws.Cells.Item($r, 12).Value = $row[2]   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $row[3]   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $row[4]   # N - Precio mínimo
    $ws.Cells.Item($r, 15).Value = $row[5]   # O - Precio máximo
    $ws.Cells.Item($r, 16).Value = $row[6]   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row[7]   # Q - Unidad de comercialización
    $ws.Cells.Item($r, 18).Value = $row[8]   # R - Origen
    $ws.Cells.Item($r, 19).Value = $row[9]   # S - Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $row[10]  # T - Kg / unidad
}
